# "final N-policy linking plots"
# Add two additional "Iteration" blocks (Iteration_1 in E:G, Iteration_2 in H:J),
# mirroring the existing Standalone block (B:D) header/interval rows, and
# refresh the underlying production-share numbers for rows 4-16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers: new merged groups for Iteration_1 (E1:G1) and Iteration_2 (H1:J1) ---
# Merge first, THEN copy formatting from an unmerged sibling cell (C1) - doing
# it in this order reuses the existing style index (no per-edge border split).
$ws.Range("E1:G1").Merge()
$ws.Range("H1:J1").Merge()

$ws.Range("E1").Value = "Iteration_1"
$ws.Range("C1").Copy()
$ws.Range("E1:G1").PasteSpecial(-4122)

$ws.Range("H1").Value = "Iteration_2"
$ws.Range("C1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)

# --- Row 2 intervals: repeat 2030/2040/2050 under each new group ---
# (text, matching the existing B2:D2 interval cells, which are strings too).
# Force text storage via NumberFormat "@" (otherwise "2030" auto-coerces to a
# number), then re-apply the B2:D2 look via PasteSpecial so the cells end up
# on the same style index as their Standalone counterparts.
$ws.Range("E2:J2").NumberFormat = "@"
$ws.Range("E2").Value = "2030"
$ws.Range("F2").Value = "2040"
$ws.Range("G2").Value = "2050"
$ws.Range("H2").Value = "2030"
$ws.Range("I2").Value = "2040"
$ws.Range("J2").Value = "2050"

$ws.Range("B2:D2").Copy()
$ws.Range("E2:G2").PasteSpecial(-4122)
$ws.Range("B2:D2").Copy()
$ws.Range("H2:J2").PasteSpecial(-4122)

# --- Refreshed production-share values (Standalone columns B:D) ---
$ws.Range("B4").Value = 1718092.7301732

$ws.Range("C6").Value = 1017699.869006407
$ws.Range("D6").Value = 1052661.907550473

$ws.Range("C8").Value = 699821.2094849477
$ws.Range("D8").Value = 664866.4223266498

# --- New Iteration_1 (E:G) and Iteration_2 (H:J) values, rows 4-16 ---
foreach ($r in 4..16) {
  foreach ($col in @("E", "F", "G", "H", "I", "J")) {
    $ws.Range("$col$r").Value = 0
  }
}

$ws.Range("E6").Value = 1717666.256704808
$ws.Range("F6").Value = 1291022.02527089
$ws.Range("G6").Value = 1699538.147168455
$ws.Range("H6").Value = 1717666.256704669
$ws.Range("I6").Value = 1291141.465709535
$ws.Range("J6").Value = 1700960.080555822

$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 426555.742260271
$ws.Range("G8").Value = 18124.34963239434
$ws.Range("H8").Value = -0.00000001909185926037406
$ws.Range("I8").Value = 426436.3265944614
$ws.Range("J8").Value = 16702.71116469113
